# Updates the "cryptos" price list worksheet to reflect the latest scrape:
#  - Column D ("Price"): refreshed price quotes for the affected coins.
#  - Column G ("Hora"): the scrape hour moves from "12" to "13" for every data row.
#
# NumberFormat is forced to Text ("@") before each assignment so that Excel stores
# the values as literal strings (preserving formatting such as trailing zeros, e.g.
# "0.9900" or "0.05950") instead of silently converting them to numeric cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: "Price" updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "251.33"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.82"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.940"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05950"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.569"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.412"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.322"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7949"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1487"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07865"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03029"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09261"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.570"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001678"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04760"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006069"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006208"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005688"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001066"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.683"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.210"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3305"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1255"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006474"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04435"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007011"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1069"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003149"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01035"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.002460"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005886"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9900"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1036"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"

# --- Column G: "Hora" updates (every row, 2-51, "12" -> "13") ---
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "13"

